# Apply the "jetzt passen die namen zusammen" edit:
# Append 32 new data rows (rows 34-65) to sheet1, with columns A (Run ID),
# B (number), C (mpn_S), and a single new G34 cell containing the new
# shared string "dienstagmorgen" (continuing the Zeit/column-G pattern that
# already has G2 = "montagnacht").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 34-65: columns A, B, C
$data = @(
    @(1,1,-3),
    @(1,2,-3),
    @(1,3,-3),
    @(1,4,-3),
    @(1,5,-3),
    @(1,6,-3),
    @(1,7,-3),
    @(1,8,-3),
    @(2,9,-3),
    @(2,10,-3),
    @(2,11,-3),
    @(2,12,-3),
    @(2,13,-3),
    @(2,14,-3),
    @(2,15,-3),
    @(2,16,-3),
    @(3,17,-3),
    @(3,18,-3),
    @(3,19,-3),
    @(3,20,-3),
    @(3,21,-3),
    @(3,22,-3),
    @(3,23,-3),
    @(3,24,-3),
    @(4,25,-3),
    @(4,26,-3),
    @(4,27,-3),
    @(4,28,-3),
    @(4,29,-3),
    @(4,30,-3),
    @(4,31,-3),
    @(4,32,-3)
)

$startRow = 34
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# New string value in column G for the first new row, matching the existing
# pattern (G1 = "Zeit" header, G2 = "montagnacht", new one = "dienstagmorgen")
$ws.Range("G34").Value = "dienstagmorgen"

# Update the view to match the author's final selection/scroll position
# (best-effort: moves the viewport and selects A66, matching the saved
# worksheet's sheetView/selection in the target workbook)
$excel.ActiveWindow.ScrollRow = 35
$ws.Range("A66").Select()
